# Stock App Project Planning Notes - "fetch watchlist detail actions"
#
# 1. Remove the empty bullet paragraph that sat just above the
#    "Resolved: ..." bullet (merges it away).
# 2. Add a new run " (separate API call)" right after the
#    "Load latest stock info in company view" bullet text.
# 3-5. The w:lastRenderedPageBreak marker (a render artifact, not part of
#    the normal Word object model) shifts from the end of one bullet to
#    the start of the next bullet in three places. Since there is no OM
#    property for it, we replace each affected paragraph's content range
#    with the equivalent OOXML via Range.InsertXML (whole-paragraph
#    range so the surrounding run rsids / pPr stay untouched).

$d = $word.ActiveDocument

function Find-ParaIndexStartsWith($doc, $prefix) {
    $i = 0
    foreach ($p in $doc.Paragraphs) {
        $i++
        $t = $p.Range.Text.TrimEnd("`r")
        if ($t.StartsWith($prefix)) {
            return $i
        }
    }
    return -1
}

function Set-ParaInnerXml($doc, $paraIndex, $bodyInnerXml) {
    # Replace a paragraph's content (everything except its own paragraph
    # mark) with literal OOXML. Targeting the *whole* paragraph range
    # (not a sub-range) keeps InsertXML's "replace in place" semantics
    # from relocating the new content to the end of the paragraph.
    $p = $doc.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $full = $doc.Range($r.Start, $r.End - 1)
    $pkg = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyInnerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $full.InsertXML($pkg)
}

# ---------------------------------------------------------------------
# Hunk 1: delete the empty ListParagraph bullet directly above "Resolved:"
# ---------------------------------------------------------------------
$resolvedIdx = Find-ParaIndexStartsWith $d "Resolved:"
$emptyIdx = $resolvedIdx - 1
$emptyPara = $d.Paragraphs.Item($emptyIdx)
if ($emptyPara.Range.Text.TrimEnd("`r").Length -eq 0) {
    $emptyPara.Range.Delete()
}

# ---------------------------------------------------------------------
# Hunk 2: "Load latest stock info in company view" gains a new trailing
# run " (separate API call)". Rebuilding the paragraph's content range
# (rather than a plain InsertAfter, which coalesces into the existing
# run, or a split/merge, which clobbers the <w:p> attributes) keeps the
# new text in its own w:r while leaving the <w:p> tag untouched.
# ---------------------------------------------------------------------
$loadIdx = Find-ParaIndexStartsWith $d "Load latest stock info in company view"
Set-ParaInnerXml $d $loadIdx '<w:p><w:r><w:t>Load latest stock info in company view</w:t></w:r><w:r><w:t xml:space="preserve"> (separate API call)</w:t></w:r></w:p>'

# ---------------------------------------------------------------------
# Hunk 3: move w:lastRenderedPageBreak from "New user form / route
# (/signup)" to the "X " run that starts the next bullet
# ---------------------------------------------------------------------
$newUserIdx = Find-ParaIndexStartsWith $d "New user form / route (/signup)"
Set-ParaInnerXml $d $newUserIdx '<w:p><w:r><w:t>New user form / route (/signup)</w:t></w:r></w:p>'

$downloadIdx = $newUserIdx + 1
Set-ParaInnerXml $d $downloadIdx '<w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">X </w:t></w:r><w:r w:rsidR="0023692C"><w:t xml:space="preserve">Started: </w:t></w:r><w:r w:rsidR="00A9214B"><w:t>Download company data and load to database / view</w:t></w:r></w:p>'

# ---------------------------------------------------------------------
# Hunk 4: move w:lastRenderedPageBreak from "3 month avg trading volume"
# to "52 week high"
# ---------------------------------------------------------------------
$threeMonthIdx = Find-ParaIndexStartsWith $d "3 month avg trading volume"
Set-ParaInnerXml $d $threeMonthIdx '<w:p><w:r><w:t>3 month avg trading volume</w:t></w:r></w:p>'

$weekHighIdx = $threeMonthIdx + 1
Set-ParaInnerXml $d $weekHighIdx '<w:p><w:r><w:lastRenderedPageBreak/><w:t>52 week high</w:t></w:r></w:p>'

# ---------------------------------------------------------------------
# Hunk 5: move w:lastRenderedPageBreak from "Chart: Displays chart" to
# "Statistics: Displays company statistics"
# ---------------------------------------------------------------------
$i = 0
$chartIdx = -1
foreach ($p in $d.Paragraphs) {
    $i++
    $t = $p.Range.Text.TrimEnd("`r")
    if ($t -eq "Chart: Displays chart") {
        $chartIdx = $i
        break
    }
}
Set-ParaInnerXml $d $chartIdx '<w:p><w:r><w:t>Chart:</w:t></w:r><w:r w:rsidR="00481AFE"><w:t xml:space="preserve"> Displays chart</w:t></w:r></w:p>'

$statsIdx = $chartIdx + 1
Set-ParaInnerXml $d $statsIdx '<w:p><w:r><w:lastRenderedPageBreak/><w:t>Statistics: Displays company statistics</w:t></w:r></w:p>'

Write-Output "done"
